$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.882.66"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.385.42"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'582.81"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'180.60"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "'0.200"
$ws.Range("E9").Value = "  +8.94%  "
$ws.Range("D10").Value = "'0.592"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").Value = "'48.68"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "'0.0000287"
$ws.Range("E12").Value = "  +4.70%  "
$ws.Range("D13").Value = "'685.07"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'8.64"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "3.933.55"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "69.884.76"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "3.388.27"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").Value = "'17.71"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "'11.34"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "'0.915"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "'17.39"
$ws.Range("E22").Value = "  +2.47%  "
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'102.13"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'2.71"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").Value = "'9.84"
$ws.Range("E27").Value = "  +3.64%  "
$ws.Range("D28").Value = "'33.75"
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("D29").Value = "'8.79"
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").Value = "'6.95"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'3.86"
$ws.Range("E31").Value = "  +15.55%  "
$ws.Range("D32").Value = "'11.12"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'558.46"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'58.18"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "3.628.15"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("D39").Value = "'35.48"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("E40").Value = "  +8.93%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.77"
$ws.Range("E41").Value = "  +5.75%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.33"
$ws.Range("E42").Value = "  +4.68%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'3.37"
$ws.Range("E43").Value = "  +3.43%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0429"
$ws.Range("E44").Value = "  +4.19%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.339"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'130.64"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'2.62"
